# Add the new "Employee" worksheet as the last sheet in the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Employee"

# Header row
$ws.Range("A1").Value = "Full Name"
$ws.Range("B1").Value = "Father Name"
$ws.Range("C1").Value = "year"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "state"

# Data row - written in E,D,B,A,C order so the shared-string table is
# appended in the same order the source workbook used.
$ws.Range("E2").Value = "Telangana"
$ws.Range("D2").Value = "T-13 Nandavanam colony"
$ws.Range("B2").Value = "Ravinder"
$ws.Range("A2").Value = "Saikumar Verramalla"
$ws.Range("C2").Value = 2025

# Column widths matching the authored sheet
$ws.Columns.Item(1).ColumnWidth = 25.7265625
$ws.Columns.Item(2).ColumnWidth = 14.453125
$ws.Columns.Item(4).ColumnWidth = 34.90625
$ws.Columns.Item(5).ColumnWidth = 12.1796875

# Selection / active cell on the new sheet
[void]$ws.Range("F17").Select()
